$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Selection / active cell (cosmetic view state) ---
$ws.Range("B32").Select()

# --- E column: clear the "No Card draw message?" notes, replace some with the
#     new "If DC active ... Restore 1 of each point ..." note, and mark a couple
#     of rows with the new blank/"checked off" fill style. ---

# Row 2: remove note, leave the existing light-grey fill as-is.
$ws.Range("E2").ClearContents()

# Row 3: remove note AND mark this row as checked-off (new plain white fill).
$ws.Range("E3").ClearContents()
$ws.Range("E3").Interior.ThemeColor = 2
$ws.Range("E3").Interior.TintAndShade = 0

# Row 4: replace the placeholder note with the real note text.
$ws.Range("E4").Value = "If DC active you can target an allied facility if not then target one of your facilities. Restore 1 of each point on the targeted facility. "

# Row 5: same note text AND mark checked-off (new plain white fill).
$ws.Range("E5").Value = "If DC active you can target an allied facility if not then target one of your facilities. Restore 1 of each point on the targeted facility. "
$ws.Range("E5").Interior.ThemeColor = 2
$ws.Range("E5").Interior.TintAndShade = 0

# --- Column B/C: check off additional cards (No/Partially -> Fully) ---
$ws.Range("C7").Value = "Fully"
$ws.Range("C8").Value = "Fully"
$ws.Range("C9").Value = "Fully"
$ws.Range("C11").Value = "Fully"
$ws.Range("C12").Value = "Fully"
$ws.Range("C13").Value = "Fully"

$ws.Range("B18").Value = "Fully"
$ws.Range("C18").Value = "Fully"

$ws.Range("B24").Value = "Fully"
$ws.Range("C24").Value = "Fully"

$ws.Range("C25").Value = "Fully"

$ws.Range("B26").Value = "Fully"
$ws.Range("C26").Value = "Fully"

# --- Row 21: remove the stray note entirely (removes the cell, not just blanks it) ---
$ws.Range("E21").ClearContents()

# --- Rows 24/26/28: remove the "Assumed because other card works" notes ---
$ws.Range("E24").ClearContents()
$ws.Range("E26").ClearContents()
$ws.Range("E28").ClearContents()
